$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces this numeric-looking value to stay text (matches
# the source file, where this cell is an inline string, not a number).
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 29.03.2024"

$ws.Range("B6").Value = "02.04."
$ws.Range("C6").Value = "03.04."
$ws.Range("D6").Value = "KARTENZ./02.04 EDEKA RO"
$ws.Range("E6").Value = "53,98-"

$ws.Range("B7").Value = "05.04."
$ws.Range("C7").Value = "06.04."
$ws.Range("D7").Value = "MCDONALDS Goslar"
$ws.Range("E7").Value = "35,40-"

$ws.Range("B8").Value = "07.04."
$ws.Range("C8").Value = "08.04."
$ws.Range("D8").Value = "BURGER KING Hildesheim"
$ws.Range("E8").Value = "42,92-"

$ws.Range("B9").Value = "09.04."
$ws.Range("C9").Value = "10.04."
$ws.Range("D9").Value = "PAYPAL JRTFJM"
$ws.Range("E9").Value = "37,26-"

$ws.Range("B10").Value = "12.04."
$ws.Range("C10").Value = "13.04."
$ws.Range("D10").Value = "KARTENZ./12.04 REWE RO"
$ws.Range("E10").Value = "128,32-"

$ws.Range("B11").Value = "16.04."
$ws.Range("C11").Value = "17.04."
$ws.Range("D11").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E11").Value = "25,39-"

$ws.Range("D12").Value = "KONTOSTAND AM 20.04.2024"
$ws.Range("E12").Value = "323,27-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 28.04.2024"
